$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while preserving it as text.
# Excel (and this COM shim) auto-converts plain-decimal-looking strings
# (e.g. "224.25", "1.00") into numbers when assigned via .Value. The source
# cells are inline/shared strings (t="inlineStr"), so we force text by using
# a leading apostrophe (Excel's standard 'treat as text' input convention)
# whenever the new text would otherwise be parsed as a plain number.
function Set-TextValue($range, [string]$text) {
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range('D2') '34.199.09'
Set-TextValue $ws.Range('E2') '  -0.95%  '
Set-TextValue $ws.Range('D3') '1.784.58'
Set-TextValue $ws.Range('E3') '  -2.93%  '
Set-TextValue $ws.Range('E4') '  +0.09%  '
Set-TextValue $ws.Range('D5') '224.25'
Set-TextValue $ws.Range('E5') '  -2.93%  '
Set-TextValue $ws.Range('D6') '0.551'
Set-TextValue $ws.Range('E6') '  +0.13%  '
Set-TextValue $ws.Range('E7') '  +0.14%  '
Set-TextValue $ws.Range('D8') '32.22'
Set-TextValue $ws.Range('E8') '  +1.50%  '
Set-TextValue $ws.Range('E9') '  -2.15%  '
Set-TextValue $ws.Range('D10') '0.0655'
Set-TextValue $ws.Range('E10') '  -3.23%  '
Set-TextValue $ws.Range('E11') '  -0.06%  '
Set-TextValue $ws.Range('D12') '2.040.26'
Set-TextValue $ws.Range('E12') '  -2.64%  '
Set-TextValue $ws.Range('D13') '11.15'
Set-TextValue $ws.Range('E13') '  +6.35%  '
Set-TextValue $ws.Range('D14') '1.797.06'
Set-TextValue $ws.Range('E14') '  -2.55%  '
Set-TextValue $ws.Range('D15') '0.625'
Set-TextValue $ws.Range('E15') '  -4.23%  '
Set-TextValue $ws.Range('D16') '34.199.06'
Set-TextValue $ws.Range('E16') '  -0.84%  '
Set-TextValue $ws.Range('E17') '  -2.09%  '
Set-TextValue $ws.Range('D18') '68.71'
Set-TextValue $ws.Range('E18') '  -2.22%  '
Set-TextValue $ws.Range('D19') '254.10'
Set-TextValue $ws.Range('E19') '  -2.99%  '
Set-TextValue $ws.Range('E20') '  -2.37%  '
Set-TextValue $ws.Range('D21') '1.00'
Set-TextValue $ws.Range('E21') '  +0.04%  '
Set-TextValue $ws.Range('E22') '  -2.71%  '
Set-TextValue $ws.Range('E23') '  -5.08%  '
Set-TextValue $ws.Range('D24') '2.13'
Set-TextValue $ws.Range('E24') '  -4.29%  '
Set-TextValue $ws.Range('D25') '157.10'
Set-TextValue $ws.Range('E25') '  -1.01%  '
Set-TextValue $ws.Range('E26') '  -2.70%  '
Set-TextValue $ws.Range('D27') '6.98'
Set-TextValue $ws.Range('E27') '  -2.51%  '
Set-TextValue $ws.Range('E28') '  -1.95%  '
Set-TextValue $ws.Range('E29') '  +0.07%  '
Set-TextValue $ws.Range('E30') '  -3.55%  '
Set-TextValue $ws.Range('E31') '  -1.85%  '
Set-TextValue $ws.Range('E32') '  -2.59%  '
Set-TextValue $ws.Range('E33') '  -1.03%  '
Set-TextValue $ws.Range('D34') '1.86'
Set-TextValue $ws.Range('E34') '  +2.85%  '
Set-TextValue $ws.Range('D35') '1.438.44'
Set-TextValue $ws.Range('E35') '  -7.38%  '
Set-TextValue $ws.Range('E36') '  -2.59%  '
Set-TextValue $ws.Range('E37') '  -1.98%  '
Set-TextValue $ws.Range('D38') '0.622'
Set-TextValue $ws.Range('E38') '  -4.15%  '
Set-TextValue $ws.Range('E39') '  +1.23%  '
Set-TextValue $ws.Range('D40') '82.76'
Set-TextValue $ws.Range('E40') '  -3.08%  '
Set-TextValue $ws.Range('E41') '  +0.44%  '
Set-TextValue $ws.Range('D42') '0.888'
Set-TextValue $ws.Range('E42') '  -3.43%  '
Set-TextValue $ws.Range('E43') '  -3.84%  '
Set-TextValue $ws.Range('E44') '  -4.02%  '
Set-TextValue $ws.Range('E45') '  -2.79%  '
Set-TextValue $ws.Range('D46') '1.940.41'
Set-TextValue $ws.Range('E46') '  -2.29%  '
Set-TextValue $ws.Range('E47') '  -0.43%  '
Set-TextValue $ws.Range('D48') '12.19'
Set-TextValue $ws.Range('E48') '  -3.11%  '
Set-TextValue $ws.Range('D49') '1.00'
Set-TextValue $ws.Range('E49') '  +0.04%  '
Set-TextValue $ws.Range('D50') '98.75'
Set-TextValue $ws.Range('E50') '  +0.28%  '
Set-TextValue $ws.Range('D51') '49.24'
Set-TextValue $ws.Range('E51') '  -6.96%  '
